$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Napoli"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "15"
$ws.Range("H2").Value = "44:14"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "47"
$ws.Range("K2").Value = "W"
$ws.Range("L2").Value = "P"
$ws.Range("C3").Value = "AC Milan"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "11"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "5"
$ws.Range("H3").Value = "35:20"
$ws.Range("J3").Value = "R"
$ws.Range("K3").Value = "R"
$ws.Range("M3").Value = "W"
$ws.Range("C4").Value = "Inter"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "18"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "12"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "1"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "5"
$ws.Range("H4").Value = "38:24"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "37"
$ws.Range("K4").Value = "R"
$ws.Range("N4").Value = "W"
$ws.Range("C5").Value = "Lazio"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "18"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "10"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "4"
$ws.Range("H5").Value = "31:15"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "34"
$ws.Range("J5").Value = "W"
$ws.Range("K5").Value = "R"
$ws.Range("L5").Value = "P"
$ws.Range("N5").Value = "W"
$ws.Range("C6").Value = "Atalanta"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "10"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "4"
$ws.Range("H6").Value = "34:20"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "34"
$ws.Range("J6").Value = "W"
$ws.Range("L6").Value = "R"
$ws.Range("M6").Value = "P"
$ws.Range("C7").Value = "AS Roma"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "18"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "10"
$ws.Range("H7").Value = "23:16"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "34"
$ws.Range("L7").Value = "W"
$ws.Range("N7").Value = "R"
$ws.Range("C8").Value = "Udinese"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "18"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "7"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "5"
$ws.Range("H8").Value = "26:21"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "25"
$ws.Range("J8").Value = "P"
$ws.Range("K8").Value = "P"
$ws.Range("L8").Value = "R"
$ws.Range("N8").Value = "R"
$ws.Range("C9").Value = "Torino"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "18"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "6"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "7"
$ws.Range("H9").Value = "18:20"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "23"
$ws.Range("M9").Value = "R"
$ws.Range("N9").Value = "W"
$ws.Range("C10").Value = "Fiorentina"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "7"
$ws.Range("H10").Value = "21:24"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "23"
$ws.Range("J10").Value = "P"
$ws.Range("K10").Value = "W"
$ws.Range("M10").Value = "P"
$ws.Range("C11").Value = "Juventus"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "11"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "3"
$ws.Range("H11").Value = "27:12"
$ws.Range("L11").Value = "W"
$ws.Range("C12").Value = "Bologna"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "4"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "8"
$ws.Range("H12").Value = "23:29"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "22"
$ws.Range("K12").Value = "P"
$ws.Range("L12").Value = "P"
$ws.Range("N12").Value = "P"
$ws.Range("C13").Value = "Empoli"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "7"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "6"
$ws.Range("H13").Value = "16:22"
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "22"
$ws.Range("J13").Value = "W"
$ws.Range("K13").Value = "R"
$ws.Range("L13").Value = "R"
$ws.Range("M13").Value = "W"
$ws.Range("C14").Value = "Monza"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "6"
$ws.Range("H14").Value = "22:27"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "21"
$ws.Range("J14").Value = "W"
$ws.Range("K14").Value = "R"
$ws.Range("C15").Value = "Lecce"
$ws.Range("D15").Value = "0:0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "19"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "4"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "9"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "6"
$ws.Range("I15").Value = "18:20"
$ws.Range("J15").Value = "?"
$ws.Range("M15").Value = "W"
$ws.Range("N15").Value = "W"
$ws.Range("C16").Value = "Spezia"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "6"
$ws.Range("H16").Value = "17:28"
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "18"
$ws.Range("J16").Value = "W"
$ws.Range("K16").Value = "R"
$ws.Range("L16").Value = "R"
$ws.Range("M16").Value = "W"
$ws.Range("C17").Value = "Salernitana"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "6"
$ws.Range("H17").Value = "23:35"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "18"
$ws.Range("J17").Value = "P"
$ws.Range("K17").Value = "R"
$ws.Range("L17").Value = "P"
$ws.Range("M17").Value = "P"
$ws.Range("C18").Value = "Sassuolo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "4"
$ws.Range("H18").Value = "17:28"
$ws.Range("I18").NumberFormat = "@"
$ws.Range("I18").Value = "16"
$ws.Range("N18").Value = "R"
$ws.Range("C19").Value = "Verona"
$ws.Range("D19").Value = "0:0"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "19"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "2"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "4"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "13"
$ws.Range("I19").Value = "15:31"
$ws.Range("J19").Value = "?"
$ws.Range("K19").Value = "P"
$ws.Range("L19").Value = "W"
$ws.Range("M19").Value = "R"
$ws.Range("N19").Value = "P"
$ws.Range("C20").Value = "Sampdoria"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "3"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "13"
$ws.Range("H20").Value = "8:31"
$ws.Range("I20").NumberFormat = "@"
$ws.Range("I20").Value = "9"
$ws.Range("K20").Value = "P"
$ws.Range("L20").Value = "W"
$ws.Range("N20").Value = "P"
$ws.Range("C21").Value = "Cremonese"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "7"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "11"
$ws.Range("H21").Value = "13:32"
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = "7"
$ws.Range("J21").Value = "P"
$ws.Range("N21").Value = "R"
